$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the two new columns, matching the style of the existing
# header cells (H1 etc.)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for column I (I0) and column J (IF) for rows 2..73
$iValues = @(10,9,9,9,9,9,8,9,9,9,9,9,9,8,8,9,9,9,9,9,9,9,7,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,10,10,9,9,9,9,9,6,9,9,9,9,11,9,9,9,9,9,9,9,9,9,4,5,4,4,5,4,3,3,2)
$jValues = @(10,9,9,9,9,9,8,9,9,9,9,9,10,8,9,9,9,9,9,9,9,9,7,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,10,10,9,9,9,9,9,7,9,9,9,9,11,9,9,10,9,9,9,9,9,9,4,5,5,4,5,4,3,3,2)

for ($n = 0; $n -lt $iValues.Length; $n++) {
    $row = $n + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$n]
    $ws.Cells.Item($row, 10).Value = $jValues[$n]
}
